$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: /category/update/:id
$ws.Range("A12").Value = "/category/update/:id"
$ws.Range("B12").Value = "put"
$ws.Range("C12").Value = "change the services in category"
$ws.Range("D12").Value = "yes"
$ws.Range("E12").Value = "Services"
$ws.Range("F12").Value = "status and message"

# Row 13: /category/delete/:id
$ws.Range("A13").Value = "/category/delete/:id"
$ws.Range("B13").Value = "delete"
$ws.Range("C13").Value = "delete a category"
$ws.Range("D13").Value = "yes"
$ws.Range("F13").Value = "status and message"

# Match style of E column (wrap text) as used in other rows (style index 4)
$ws.Range("E12").WrapText = $true

$ws.Range("G13").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
